$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A289").Value = "Why can't I add 251 curve shades to my log?"
$ws.Range("B289").Value = "llama3.2:latest"
$ws.Range("C289").Value = "You cannot add 251 curve shades to your log because the limit is 250."

$ws.Range("A290").Value = "I have 20000 modifiers added ty log, why I can't I add anymore?"
$ws.Range("B290").Value = "llama3.2:latest"
$ws.Range("C290").Value = "You cannot add more than 20000 modifiers per plot because of the limit on the number of layouts per ODF file. The maximum number of layouts allowed is 20000."

$ws.Range("A291").Value = "I have 20000 modifiers added ty log, why I can't I add anymore?"
$ws.Range("B291").Value = "llama3.2:latest"
$ws.Range("C291").Value = "You cannot add more than 20000 modifiers per plot because of the limit on the number of layouts per ODF file. The maximum number of layouts allowed is 20000."

$ws.Range("A292").Value = "How many log headers can I add to my log?"
$ws.Range("B292").Value = "llama3.2:latest"
$ws.Range("C292").Value = "Based on the provided feedback, you can select up to 50 headers for display in your log."
